$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in row 5 with the new "no conformidad" entry (mirrors row 4's layout)
$ws.Range("B5").Value = "No todos los tickets tienen un nombre descriptivo"
$ws.Range("C5").Value = "Oriana Osiris"
$ws.Range("D5").Value = 42366
$ws.Range("E5").Value = 42366
$ws.Range("D5").NumberFormat = "DD/MM/YY"
$ws.Range("E5").NumberFormat = "DD/MM/YY"
$ws.Range("F5").Value = "Cerrada"
$ws.Range("G5").Value = "Se genera aviso de cambiar el nombre a uno mas descriptivo"

# Row grows taller to fit the wrapped text, matching row 4's height
$ws.Rows.Item(5).RowHeight = 28.35

# Update the active selection left by the editor
$ws.Range("F6").Select()

$wb.Save()
